$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the company name in B2 (was "ABL Imaging Group Inc.") with the new value.
$ws.Range("B2").Select()
$excel.ActiveCell.Value = "XYZ Test Inc"
